$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename sheet Sheet2 -> Sheet1
$ws.Name = "Sheet1"

# 2. Insert two new (blank) columns at the very start (A:B), shifting
#    everything right by 2. This also shifts the <cols> width metadata,
#    which matches the target diff exactly.
$ws.Columns("A:B").Insert()

# At this point the header row (old B1..AM1, "lelaki biasa".."tgl buat")
# now lives in D1..AO1. Old A1 ("no tp") now lives in C1.

# 3. Set the new headers for the first three columns (do this before the
#    shared-string-producing rename below so new strings intern in the
#    same order the target workbook used: "periode" / "no tp" already
#    exist, but "no ba" is brand new and must get the lower index).
$ws.Cells.Item(1, 1).Value = "no ba"
$ws.Cells.Item(1, 2).Value = "no tp"
$ws.Cells.Item(1, 3).Value = "periode"

# 4. Insert a new column "aset_likuid_tidak_menghasilkan" right after the
#    "aset tidak menghasilkan" column (old J1, now L1), WITHOUT using a
#    structural column insert (that would incorrectly shift <cols> widths
#    that must stay put per the target diff). Instead, shift the cell
#    values manually column-by-column (back to front) then write the new
#    header text into the freed-up cell.
for ($col = 39; $col -ge 13; $col--) {
    $srcVal = $ws.Cells.Item(1, $col).Value2
    $ws.Cells.Item(1, $col + 1).Value = $srcVal
}
$ws.Cells.Item(1, 13).Value = "aset_likuid_tidak_menghasilkan"

# 5. Rename "tgl buat" -> "tanggal buat" (now at the very end, column AO = 41)
$ws.Cells.Item(1, 41).Value = "tanggal buat"

# 6. Update the selection to match the target sheetView.
$ws.Range("C18").Select()
